# Refresh the toy-spam anchor-word tables for min_count=5.
# Column layout: A:H = negative-anchor words, J:Q = positive-anchor words.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Negative word table (A3:H32) ---
$negArr = New-Object 'object[,]' 30,8
$negArr[0,0] = "poorly"
$negArr[0,1] = 0.9782608695652174
$negArr[0,2] = 45
$negArr[0,3] = 45
$negArr[0,4] = 0
$negArr[0,5] = 1
$negArr[0,6] = $false
$negArr[0,7] = 1
$negArr[1,0] = "disappointing"
$negArr[1,1] = 0.8636363636363636
$negArr[1,2] = 38
$negArr[1,3] = 38
$negArr[1,4] = 0
$negArr[1,5] = 1
$negArr[1,6] = $false
$negArr[1,7] = 6
$negArr[2,0] = "broke"
$negArr[2,1] = 0.7815533980582524
$negArr[2,2] = 161
$negArr[2,3] = 161
$negArr[2,4] = 0
$negArr[2,5] = 1
$negArr[2,6] = $false
$negArr[2,7] = 45
$negArr[3,0] = "however"
$negArr[3,1] = 0.75
$negArr[3,2] = 48
$negArr[3,3] = 48
$negArr[3,4] = 0
$negArr[3,5] = 1
$negArr[3,6] = $false
$negArr[3,7] = 16
$negArr[4,0] = "disappointed"
$negArr[4,1] = 0.7473118279569892
$negArr[4,2] = 139
$negArr[4,3] = 139
$negArr[4,4] = 0
$negArr[4,5] = 1
$negArr[4,6] = $false
$negArr[4,7] = 47
$negArr[5,0] = "poor"
$negArr[5,1] = 0.7323943661971831
$negArr[5,2] = 52
$negArr[5,3] = 52
$negArr[5,4] = 0
$negArr[5,5] = 1
$negArr[5,6] = $false
$negArr[5,7] = 19
$negArr[6,0] = "waste"
$negArr[6,1] = 0.6756756756756757
$negArr[6,2] = 100
$negArr[6,3] = 100
$negArr[6,4] = 0
$negArr[6,5] = 1
$negArr[6,6] = $false
$negArr[6,7] = 48
$negArr[7,0] = "smaller"
$negArr[7,1] = 0.6050420168067226
$negArr[7,2] = 72
$negArr[7,3] = 72
$negArr[7,4] = 0
$negArr[7,5] = 1
$negArr[7,6] = $false
$negArr[7,7] = 47
$negArr[8,0] = "instead"
$negArr[8,1] = 0.6041666666666666
$negArr[8,2] = 29
$negArr[8,3] = 29
$negArr[8,4] = 0
$negArr[8,5] = 1
$negArr[8,6] = $false
$negArr[8,7] = 19
$negArr[9,0] = "junk"
$negArr[9,1] = 0.5272727272727272
$negArr[9,2] = 29
$negArr[9,3] = 29
$negArr[9,4] = 0
$negArr[9,5] = 1
$negArr[9,6] = $false
$negArr[9,7] = 26
$negArr[10,0] = "broken"
$negArr[10,1] = 0.5060240963855421
$negArr[10,2] = 42
$negArr[10,3] = 42
$negArr[10,4] = 0
$negArr[10,5] = 1
$negArr[10,6] = $false
$negArr[10,7] = 41
$negArr[11,0] = "small"
$negArr[11,1] = 0.4985507246376812
$negArr[11,2] = 172
$negArr[11,3] = 172
$negArr[11,4] = 0
$negArr[11,5] = 1
$negArr[11,6] = $false
$negArr[11,7] = 173
$negArr[12,0] = "plastic"
$negArr[12,1] = 0.4251968503937008
$negArr[12,2] = 54
$negArr[12,3] = 54
$negArr[12,4] = 0
$negArr[12,5] = 1
$negArr[12,6] = $false
$negArr[12,7] = 73
$negArr[13,0] = "apart"
$negArr[13,1] = 0.4210526315789473
$negArr[13,2] = 40
$negArr[13,3] = 40
$negArr[13,4] = 0
$negArr[13,5] = 1
$negArr[13,6] = $false
$negArr[13,7] = 55
$negArr[14,0] = "ok"
$negArr[14,1] = 0.34375
$negArr[14,2] = 44
$negArr[14,3] = 44
$negArr[14,4] = 0
$negArr[14,5] = 1
$negArr[14,6] = $false
$negArr[14,7] = 84
$negArr[15,0] = "thought"
$negArr[15,1] = 0.3069306930693069
$negArr[15,2] = 62
$negArr[15,3] = 62
$negArr[15,4] = 0
$negArr[15,5] = 1
$negArr[15,6] = $false
$negArr[15,7] = 140
$negArr[16,0] = "cheap"
$negArr[16,1] = 0.2890995260663507
$negArr[16,2] = 61
$negArr[16,3] = 61
$negArr[16,4] = 0
$negArr[16,5] = 1
$negArr[16,6] = $false
$negArr[16,7] = 150
$negArr[17,0] = "though"
$negArr[17,1] = 0.2564102564102564
$negArr[17,2] = 30
$negArr[17,3] = 30
$negArr[17,4] = 0
$negArr[17,5] = 1
$negArr[17,6] = $false
$negArr[17,7] = 87
$negArr[18,0] = "size"
$negArr[18,1] = 0.2164948453608248
$negArr[18,2] = 42
$negArr[18,3] = 42
$negArr[18,4] = 0
$negArr[18,5] = 1
$negArr[18,6] = $false
$negArr[18,7] = 152
$negArr[19,0] = "used"
$negArr[19,1] = 0.1942857142857143
$negArr[19,2] = 34
$negArr[19,3] = 34
$negArr[19,4] = 0
$negArr[19,5] = 1
$negArr[19,6] = $false
$negArr[19,7] = 141
$negArr[20,0] = "hard"
$negArr[20,1] = 0.185
$negArr[20,2] = 37
$negArr[20,3] = 37
$negArr[20,4] = 0
$negArr[20,5] = 1
$negArr[20,6] = $false
$negArr[20,7] = 163
$negArr[21,0] = "item"
$negArr[21,1] = 0.1847826086956522
$negArr[21,2] = 51
$negArr[21,3] = 51
$negArr[21,4] = 0
$negArr[21,5] = 1
$negArr[21,6] = $false
$negArr[21,7] = 225
$negArr[22,0] = "would"
$negArr[22,1] = 0.1795252225519288
$negArr[22,2] = 121
$negArr[22,3] = 121
$negArr[22,4] = 0
$negArr[22,5] = 1
$negArr[22,6] = $false
$negArr[22,7] = 553
$negArr[23,0] = "money"
$negArr[23,1] = 0.1487341772151899
$negArr[23,2] = 47
$negArr[23,3] = 47
$negArr[23,4] = 0
$negArr[23,5] = 1
$negArr[23,6] = $false
$negArr[23,7] = 269
$negArr[24,0] = "work"
$negArr[24,1] = 0.1455696202531646
$negArr[24,2] = 46
$negArr[24,3] = 46
$negArr[24,4] = 0
$negArr[24,5] = 1
$negArr[24,6] = $false
$negArr[24,7] = 270
$negArr[25,0] = "better"
$negArr[25,1] = 0.1355140186915888
$negArr[25,2] = 29
$negArr[25,3] = 29
$negArr[25,4] = 0
$negArr[25,5] = 1
$negArr[25,6] = $false
$negArr[25,7] = 185
$negArr[26,0] = "product"
$negArr[26,1] = 0.13215859030837
$negArr[26,2] = 60
$negArr[26,3] = 60
$negArr[26,4] = 0
$negArr[26,5] = 1
$negArr[26,6] = $false
$negArr[26,7] = 394
$negArr[27,0] = "price"
$negArr[27,1] = 0.1235632183908046
$negArr[27,2] = 43
$negArr[27,3] = 43
$negArr[27,4] = 0
$negArr[27,5] = 1
$negArr[27,6] = $false
$negArr[27,7] = 305
$negArr[28,0] = "use"
$negArr[28,1] = 0.09863013698630137
$negArr[28,2] = 36
$negArr[28,3] = 36
$negArr[28,4] = 0
$negArr[28,5] = 1
$negArr[28,6] = $false
$negArr[28,7] = 329
$negArr[29,0] = "like"
$negArr[29,1] = 0.06589785831960461
$negArr[29,2] = 40
$negArr[29,3] = 41
$negArr[29,4] = 0.02
$negArr[29,5] = 0.98
$negArr[29,6] = $true
$negArr[29,7] = 567
$ws.Range("A3:H32").Value = $negArr

# --- Positive word table (J3:Q16) ---
# The positive table grows from 13 to 14 data rows, so row 16 (J:Q) needs the
# same header-style formatting (bold, thin border, centered) already used by
# J3:J15. Copy that formatting down from row 15 before writing the new values.
$ws.Range("J15:Q15").Copy() | Out-Null
$ws.Range("J16:Q16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$posArr = New-Object 'object[,]' 14,8
$posArr[0,0] = "awesome"
$posArr[0,1] = 0.8461538461538461
$posArr[0,2] = 55
$posArr[0,3] = 55
$posArr[0,4] = 1
$posArr[0,5] = 0
$posArr[0,6] = $false
$posArr[0,7] = 10
$posArr[1,0] = "wonderful"
$posArr[1,1] = 0.8214285714285714
$posArr[1,2] = 46
$posArr[1,3] = 46
$posArr[1,4] = 1
$posArr[1,5] = 0
$posArr[1,6] = $false
$posArr[1,7] = 10
$posArr[2,0] = "favorite"
$posArr[2,1] = 0.6881720430107527
$posArr[2,2] = 64
$posArr[2,3] = 64
$posArr[2,4] = 1
$posArr[2,5] = 0
$posArr[2,6] = $false
$posArr[2,7] = 29
$posArr[3,0] = "classic"
$posArr[3,1] = 0.6415094339622641
$posArr[3,2] = 34
$posArr[3,3] = 34
$posArr[3,4] = 1
$posArr[3,5] = 0
$posArr[3,6] = $false
$posArr[3,7] = 19
$posArr[4,0] = "excellent"
$posArr[4,1] = 0.5
$posArr[4,2] = 32
$posArr[4,3] = 32
$posArr[4,4] = 1
$posArr[4,5] = 0
$posArr[4,6] = $false
$posArr[4,7] = 32
$posArr[5,0] = "thank"
$posArr[5,1] = 0.463768115942029
$posArr[5,2] = 32
$posArr[5,3] = 32
$posArr[5,4] = 1
$posArr[5,5] = 0
$posArr[5,6] = $false
$posArr[5,7] = 37
$posArr[6,0] = "great"
$posArr[6,1] = 0.3360655737704918
$posArr[6,2] = 410
$posArr[6,3] = 410
$posArr[6,4] = 1
$posArr[6,5] = 0
$posArr[6,6] = $false
$posArr[6,7] = 810
$posArr[7,0] = "love"
$posArr[7,1] = 0.2926829268292683
$posArr[7,2] = 204
$posArr[7,3] = 204
$posArr[7,4] = 1
$posArr[7,5] = 0
$posArr[7,6] = $false
$posArr[7,7] = 493
$posArr[8,0] = "loves"
$posArr[8,1] = 0.2593360995850623
$posArr[8,2] = 125
$posArr[8,3] = 125
$posArr[8,4] = 1
$posArr[8,5] = 0
$posArr[8,6] = $false
$posArr[8,7] = 357
$posArr[9,0] = "best"
$posArr[9,1] = 0.2583333333333334
$posArr[9,2] = 31
$posArr[9,3] = 31
$posArr[9,4] = 1
$posArr[9,5] = 0
$posArr[9,6] = $false
$posArr[9,7] = 89
$posArr[10,0] = "perfect"
$posArr[10,1] = 0.1987951807228916
$posArr[10,2] = 33
$posArr[10,3] = 33
$posArr[10,4] = 1
$posArr[10,5] = 0
$posArr[10,6] = $false
$posArr[10,7] = 133
$posArr[11,0] = "loved"
$posArr[11,1] = 0.1773700305810398
$posArr[11,2] = 58
$posArr[11,3] = 58
$posArr[11,4] = 1
$posArr[11,5] = 0
$posArr[11,6] = $false
$posArr[11,7] = 269
$posArr[12,0] = "fun"
$posArr[12,1] = 0.08326029798422437
$posArr[12,2] = 95
$posArr[12,3] = 95
$posArr[12,4] = 1
$posArr[12,5] = 0
$posArr[12,6] = $false
$posArr[12,7] = 1046
$posArr[13,0] = "game"
$posArr[13,1] = 0.03636363636363636
$posArr[13,2] = 56
$posArr[13,3] = 57
$posArr[13,4] = 0.98
$posArr[13,5] = 0.02000000000000002
$posArr[13,6] = $true
$posArr[13,7] = 1484
$ws.Range("J3:Q16").Value = $posArr
